$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to the target cells so that values are preserved exactly as text
# (matching the inlineStr representation of the source data) instead of being converted to numbers.

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "371"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "829537.07"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "787"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2207018.22"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "498"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1291475.69"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "101"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "227455.66"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "243"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "676181.77"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "97"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "238245.00"

$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "131"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "306000.00"

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "276"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "791752.33"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "139"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350900.26"

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "6"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12700.00"

$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "178"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "395800.00"

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "289"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "741100.74"

$ws.Range("C65").NumberFormat = "@"
$ws.Range("C65").Value = "740"
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = "2136729.03"

$ws.Range("C66").NumberFormat = "@"
$ws.Range("C66").Value = "429"
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "1172472.79"

$ws.Range("C75").NumberFormat = "@"
$ws.Range("C75").Value = "152"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = "331000.00"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "402"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "1099396.01"

$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "153"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "367677.09"
